# Update the division-problem answers in the single table on the page.
# The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17 contain
# text (the other rows are blank spacer rows). Each text cell is replaced
# positionally (row, column) with its new value, since several of the old
# values are duplicated elsewhere in the table and a plain text Find/Replace
# could not disambiguate them reliably.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "89÷5=17, 4"
$t.Cell(1,2).Range.Text = "62÷6=10, 2"
$t.Cell(1,3).Range.Text = "55÷3=18, 1"
$t.Cell(1,4).Range.Text = "60÷2=30, 0"
$t.Cell(1,5).Range.Text = "19÷4=4, 3"

$t.Cell(5,1).Range.Text = "91÷4=22, 3"
$t.Cell(5,2).Range.Text = "74÷3=24, 2"
$t.Cell(5,3).Range.Text = "30÷7=4, 2"
$t.Cell(5,4).Range.Text = "56÷8=7, 0"
$t.Cell(5,5).Range.Text = "42÷9=4, 6"

$t.Cell(9,1).Range.Text = "55÷2=27, 1"
$t.Cell(9,2).Range.Text = "87÷4=21, 3"
$t.Cell(9,3).Range.Text = "63÷6=10, 3"
$t.Cell(9,4).Range.Text = "18÷5=3, 3"
$t.Cell(9,5).Range.Text = "31÷3=10, 1"

$t.Cell(13,1).Range.Text = "31÷8=3, 7"
$t.Cell(13,2).Range.Text = "76÷9=8, 4"
$t.Cell(13,3).Range.Text = "69÷6=11, 3"
$t.Cell(13,4).Range.Text = "96÷9=10, 6"
$t.Cell(13,5).Range.Text = "83÷9=9, 2"

$t.Cell(17,1).Range.Text = "61÷4=15, 1"
$t.Cell(17,2).Range.Text = "35÷9=3, 8"
$t.Cell(17,3).Range.Text = "24÷9=2, 6"
$t.Cell(17,4).Range.Text = "89÷9=9, 8"
$t.Cell(17,5).Range.Text = "53÷6=8, 5"

Write-Output "Updated 25 cells"
